$wb = $excel.ActiveWorkbook

$wsAtc      = $wb.Worksheets.Item("ATC_CDV")
$wsArticles = $wb.Worksheets.Item("Articles")
$wsAppro    = $wb.Worksheets.Item("supplier_appro")

# --- Populate Articles sheet (sheet2) with article numbers + N/R flag ---
$articlesData = @(
    @(6509975, "N"),
    @(6527018, "R"),
    @(6527118, "N"),
    @(6528319, "R"),
    @(6522478, "N"),
    @(6522578, "R"),
    @(6524481, "N"),
    @(6511396, "R"),
    @(6512197, "N")
)
for ($i = 0; $i -lt $articlesData.Length; $i++) {
    $row = $i + 1
    $wsArticles.Cells.Item($row, 1).Value = $articlesData[$i][0]
    $wsArticles.Cells.Item($row, 2).Value = $articlesData[$i][1]
}
# Reuse the existing "vertical center / wrap text" style (same as ATC_CDV!C1) for column A
$wsAtc.Range("C1").Copy()
$wsArticles.Range("A1:A9").PasteSpecial(-4122)

# --- Populate supplier_appro sheet (sheet3) with supplier codes ---
$approData = @("A3861", "B2041", "B2316", "B6185", "B6340", "B0176", "B1392", "B0553")
for ($i = 0; $i -lt $approData.Length; $i++) {
    $row = $i + 1
    $wsAppro.Cells.Item($row, 1).Value = $approData[$i]
}
$wsAtc.Range("C1").Copy()
$wsAppro.Range("A1:A8").PasteSpecial(-4122)

# --- Add new sheet "supplier_producer" as the last sheet ---
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsAppro)
$wsNew.Name = "supplier_producer"

$producerData = @("06687", "06686", "06685", "06684", "06683", "06682", "06681", "06680", "06679")
for ($i = 0; $i -lt $producerData.Length; $i++) {
    $row = $i + 1
    $wsNew.Cells.Item($row, 1).Value = $producerData[$i]
}
# New style: Text number format + vertical-center / wrap-text alignment
$wsNew.Range("A1").NumberFormat = "@"
$wsNew.Range("A1").WrapText = $true
$wsNew.Range("A1").VerticalAlignment = -4108
$wsNew.Range("A1").Copy()
$wsNew.Range("A2:A9").PasteSpecial(-4122)
# Re-set A1's value (PasteSpecial above only touched A2:A9, but make sure text stays text)
$wsNew.Range("A1").Value = "06687"

# --- Selections on each sheet (applied before the final active-sheet selection) ---
$wsAtc.Range("F4").Select()
$wsArticles.Range("E3").Select()
$wsNew.Range("A9").Select()

# supplier_appro stays the active tab with A2 selected
$wsAppro.Select()
$wsAppro.Range("A2").Select()
